# Generate Report for Handoff
#
# A second source file (ffffa1d33e43-db34-48bb-8016-75a9096b13b1.md) has
# been picked up for localization alongside the existing
# 3293cdd9-9633-45b3-a736-b0280233012c.md (renamed from
# 9dd30674-4adc-447b-b498-6616376d2697.md), and a fresh handoff round
# produced new target (.xlf) files. This inserts a row for the new file on
# every sheet (pushing ".localization-config" down one row) and updates the
# renamed/handed-off file names + timestamps.

$wb = $excel.ActiveWorkbook

function Set-CellHyperlink {
    param($ws, $cellRef, $text, $url)

    $ws.Range($cellRef).Value = $text
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $text) | Out-Null
    $ws.Range($cellRef).Font.Underline = $true
    $ws.Range($cellRef).Font.Color = 15570276
}

$commit = "14d381c84cc473408edaf8c055d132ff9c2f9a96"
$oldMd = "9dd30674-4adc-447b-b498-6616376d2697.md"
$newMd = "3293cdd9-9633-45b3-a736-b0280233012c.md"
$secondMd = "ffffa1d33e43-db34-48bb-8016-75a9096b13b1.md"
$cfgName = ".localization-config"

$zhOldXlf = "9dd30674-4adc-447b-b498-6616376d2697.99803a9f9145700dae47ef704869ab13733de735.zh-cn.xlf"
$zhNewXlf = "3293cdd9-9633-45b3-a736-b0280233012c.39e0d787f1d78e8949791b790cea43ef14b88a2f.zh-cn.xlf"
$deOldXlf = "9dd30674-4adc-447b-b498-6616376d2697.99803a9f9145700dae47ef704869ab13733de735.de-de.xlf"
$deNewXlf = "3293cdd9-9633-45b3-a736-b0280233012c.39e0d787f1d78e8949791b790cea43ef14b88a2f.de-de.xlf"

$zhTimestamp = "2016-03-09 15:53:54"
$deTimestamp = "2016-03-09 15:54:03"
$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop the existing hyperlinks on this sheet - they'll be re-created below
# pointing at the correct rows/targets (row 3 onward shifts down by one).
$wsOverview.Range("A1").Hyperlinks.Delete() | Out-Null

Set-CellHyperlink $wsOverview "A2" $newMd "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$newMd"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

Set-CellHyperlink $wsOverview "A3" $secondMd "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$secondMd"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

Set-CellHyperlink $wsOverview "A4" $cfgName "https://github.com/OpenLocalizationTest/oltest/blob/$commit/$cfgName"
$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A1").Hyperlinks.Delete() | Out-Null

Set-CellHyperlink $wsZh "A2" $newMd "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$newMd"
$wsZh.Range("B2").Value = "Ready for handoff"
Set-CellHyperlink $wsZh "C2" $zhNewXlf "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e520a2a2e82f69d032d0924b0356a545a5fabf51/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhNewXlf"
$wsZh.Range("D2").Value = $zhTimestamp
$wsZh.Range("G2").Value = $epoch
$wsZh.Range("H2").Value = "Include"

Set-CellHyperlink $wsZh "A3" $secondMd "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$secondMd"
$wsZh.Range("B3").Value = "Ready for handoff"
Set-CellHyperlink $wsZh "C3" $zhNewXlf "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e520a2a2e82f69d032d0924b0356a545a5fabf51/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhNewXlf"
$wsZh.Range("D3").Value = $zhTimestamp
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = "Include"

Set-CellHyperlink $wsZh "A4" $cfgName "https://github.com/OpenLocalizationTest/oltest/blob/$commit/$cfgName"
$wsZh.Range("B4").Value = "Not to be localized"
$wsZh.Range("D4").Value = $epoch
$wsZh.Range("G4").Value = $epoch
$wsZh.Range("H4").Value = "Ignored"

# ---------------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A1").Hyperlinks.Delete() | Out-Null

Set-CellHyperlink $wsDe "A2" $newMd "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$newMd"
$wsDe.Range("B2").Value = "Ready for handoff"
Set-CellHyperlink $wsDe "C2" $deNewXlf "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f1d0f0d32bc89157ed3959241fff9562d11b3d14/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deNewXlf"
$wsDe.Range("D2").Value = $deTimestamp
$wsDe.Range("G2").Value = $epoch
$wsDe.Range("H2").Value = "Include"

Set-CellHyperlink $wsDe "A3" $secondMd "https://github.com/OpenLocalizationTest/oltest/blob/$commit/e2e/$secondMd"
$wsDe.Range("B3").Value = "Ready for handoff"
Set-CellHyperlink $wsDe "C3" $deNewXlf "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f1d0f0d32bc89157ed3959241fff9562d11b3d14/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deNewXlf"
$wsDe.Range("D3").Value = $deTimestamp
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = "Include"

Set-CellHyperlink $wsDe "A4" $cfgName "https://github.com/OpenLocalizationTest/oltest/blob/$commit/$cfgName"
$wsDe.Range("B4").Value = "Not to be localized"
$wsDe.Range("D4").Value = $epoch
$wsDe.Range("G4").Value = $epoch
$wsDe.Range("H4").Value = "Ignored"

Write-Host "Handoff report regenerated."
